$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "added harvard case classification": recompute the *_old score columns
# (Ada_old = C, Avey_old = F, K health_old = M, WebMD_old = Q, doctor_NJ_old = U)
# for rows 2-6 (precision, recall, f1-score, f2-score, NDCG).

$ws.Range("C2").Value = 0.75
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.25
$ws.Range("Q2").Value = 0.4
$ws.Range("U2").Value = 1

$ws.Range("C3").Value = 0.75
$ws.Range("F3").Value = 1
$ws.Range("M3").Value = 0.25
$ws.Range("Q3").Value = 0.5
$ws.Range("U3").Value = 0.5

$ws.Range("C4").Value = 0.75
$ws.Range("F4").Value = 0.8
$ws.Range("M4").Value = 0.25
$ws.Range("Q4").Value = 0.4444444444444445
$ws.Range("U4").Value = 0.6666666666666666

$ws.Range("C5").Value = 0.75
$ws.Range("F5").Value = 0.9090909090909091
$ws.Range("M5").Value = 0.25
$ws.Range("Q5").Value = 0.4761904761904762
$ws.Range("U5").Value = 0.5555555555555556

$ws.Range("C6").Value = 0.2410465689186769
$ws.Range("F6").Value = 1
$ws.Range("M6").Value = 0.04684458432433119
$ws.Range("Q6").Value = 0.3981789667568151
$ws.Range("U6").Value = 0.4901792149829458
